$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '70.095.97'
$ws.Range('E2').Value = '  -0.07%  '
Set-TextValue 'D3' '3.540.84'
Set-TextValue 'D4' '0.999'
$ws.Range('E4').Value = '  -0.06%  '
Set-TextValue 'D5' '604.11'
$ws.Range('E5').Value = '  -2.11%  '
Set-TextValue 'D6' '196.89'
$ws.Range('E6').Value = '  +5.89%  '
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('E9').Value = '  -3.53%  '
Set-TextValue 'D10' '0.655'
$ws.Range('E10').Value = '  -0.37%  '
$ws.Range('E11').Value = '  +1.05%  '
Set-TextValue 'D12' '0.0000303'
$ws.Range('E12').Value = '  -1.10%  '
Set-TextValue 'D13' '9.55'
$ws.Range('E13').Value = '  -0.48%  '
Set-TextValue 'D14' '4.107.85'
$ws.Range('E14').Value = '  +0.23%  '
Set-TextValue 'D15' '603.74'
$ws.Range('E15').Value = '  -2.25%  '
Set-TextValue 'D16' '19.24'
$ws.Range('E16').Value = '  +1.34%  '
Set-TextValue 'D17' '70.206.77'
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('E18').Value = '  -0.64%  '
Set-TextValue 'D19' '3.537.14'
$ws.Range('E19').Value = '  +0.29%  '
$ws.Range('E20').Value = '  +0.59%  '
Set-TextValue 'D21' '0.995'
$ws.Range('E21').Value = '  +0.03%  '
Set-TextValue 'D22' '17.99'
$ws.Range('E22').Value = '  +2.73%  '
Set-TextValue 'D23' '5.27'
$ws.Range('E23').Value = '  +4.43%  '
Set-TextValue 'D24' '102.79'
$ws.Range('E24').Value = '  -0.75%  '
Set-TextValue 'D25' '4.62'
$ws.Range('E25').Value = '  -1.96%  '
$ws.Range('E26').Value = '  +4.09%  '
Set-TextValue 'D27' '10.98'
$ws.Range('E27').Value = '  -0.02%  '
Set-TextValue 'D28' '9.63'
$ws.Range('E28').Value = '  -1.93%  '
Set-TextValue 'D29' '33.79'
$ws.Range('E29').Value = '  -0.24%  '
Set-TextValue 'D30' '4.38'
$ws.Range('E30').Value = '  +21.93%  '
Set-TextValue 'D31' '7.14'
$ws.Range('E31').Value = '  +1.01%  '
Set-TextValue 'D32' '12.65'
$ws.Range('E32').Value = '  +2.17%  '
$ws.Range('E33').Value = '  -0.04%  '
Set-TextValue 'D34' '63.40'
$ws.Range('E34').Value = '  -1.29%  '
Set-TextValue 'D35' '0.0₃0839'
$ws.Range('E35').Value = '  +7.83%  '
Set-TextValue 'D36' '3.777.44'
$ws.Range('E36').Value = '  +6.91%  '
$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D37' '1.00'
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('B38').Value = 'Fetch.AI'
$ws.Range('C38').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D38' '3.07'
$ws.Range('E38').Value = '  -4.08%  '
Set-TextValue 'D39' '3.66'
$ws.Range('E39').Value = '  +2.20%  '
$ws.Range('E40').Value = '  -1.15%  '
Set-TextValue 'D41' '36.76'
$ws.Range('E41').Value = '  -1.38%  '
Set-TextValue 'D42' '490.94'
$ws.Range('E42').Value = '  -7.84%  '
$ws.Range('E43').Value = '  -2.48%  '
Set-TextValue 'D44' '0.0459'
$ws.Range('E44').Value = '  -1.27%  '
$ws.Range('E45').Value = '  -3.66%  '
$ws.Range('E46').Value = '  -1.78%  '
Set-TextValue 'D47' '3.30'
$ws.Range('E48').Value = '  +0.21%  '
Set-TextValue 'D49' '8.68'
$ws.Range('E49').Value = '  -4.04%  '
$ws.Range('E50').Value = '  +3.20%  '
Set-TextValue 'D51' '130.21'
$ws.Range('E51').Value = '  -2.79%  '
